$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row field updates scraped from the latest coinranking.com snapshot.
# Price (column D) cells hold free-form numeric-looking text (e.g. "64.804.27",
# "0.999", "1.00") that must stay literal text -- a bare numeric assignment would
# let Excel re-interpret/re-format it as a real number and silently mangle the
# value (e.g. "1.00" -> 1). Prefixing with an apostrophe forces text entry; the
# follow-up Style reset clears the resulting quote-prefix formatting so the cell
# style stays identical to its original (unstyled) state.
$rows = @(
    @{ Row=2; D='64.804.27'; E='  +5.23%  ' }
    @{ Row=3; D='3.095.57'; E='  +3.10%  ' }
    @{ Row=4; E='  -0.02%  ' }
    @{ Row=5; D='558.79'; E='  +3.10%  ' }
    @{ Row=6; D='143.78'; E='  +9.59%  ' }
    @{ Row=7; E='  -0.10%  ' }
    @{ Row=8; D='3.094.03'; E='  +3.27%  ' }
    @{ Row=9; D='0.499'; E='  +2.25%  ' }
    @{ Row=10; D='7.17'; E='  +18.75%  ' }
    @{ Row=11; E='  +4.98%  ' }
    @{ Row=12; D='0.463'; E='  +4.20%  ' }
    @{ Row=13; D='0.0000228'; E='  +4.42%  ' }
    @{ Row=14; D='35.29'; E='  +3.09%  ' }
    @{ Row=15; D='3.602.78'; E='  +3.24%  ' }
    @{ Row=16; D='64.796.64'; E='  +5.12%  ' }
    @{ Row=17; D='3.101.05'; E='  +3.33%  ' }
    @{ Row=18; E='  -0.57%  ' }
    @{ Row=19; D='6.80'; E='  +3.08%  ' }
    @{ Row=20; D='483.03'; E='  -0.17%  ' }
    @{ Row=21; D='13.83'; E='  +4.73%  ' }
    @{ Row=22; D='0.676'; E='  +1.48%  ' }
    @{ Row=23; D='7.55'; E='  +8.88%  ' }
    @{ Row=24; E='  +12.35%  ' }
    @{ Row=25; D='80.91'; E='  -1.51%  ' }
    @{ Row=26; D='0.999'; E='  +0.04%  ' }
    @{ Row=27; D='2.78'; E='  +3.84%  ' }
    @{ Row=28; D='8.17'; E='  +6.97%  ' }
    @{ Row=29; D='2.06'; E='  +8.42%  ' }
    @{ Row=30; D='1.00'; E='  +0.06%  ' }
    @{ Row=31; D='26.08'; E='  +1.73%  ' }
    @{ Row=32; E='  +3.75%  ' }
    @{ Row=33; D='2.46'; E='  +6.23%  ' }
    @{ Row=34; D='5.70'; E='  +2.12%  ' }
    @{ Row=35; D='6.22'; E='  +6.53%  ' }
    @{ Row=36; D='54.91'; E='  +0.24%  ' }
    @{ Row=37; D='464.73'; E='  +6.81%  ' }
    @{ Row=38; D='0.0408'; E='  +6.99%  ' }
    @{ Row=39; D='0.0824'; E='  +4.09%  ' }
    @{ Row=40; D='2.87'; E='  +19.44%  ' }
    @{ Row=41; D='3.006.94'; E='  -3.99%  ' }
    @{ Row=42; B='Cosmos'; C='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D='8.26'; E='  +2.64%  ' }
    @{ Row=43; B='Kaspa'; C='https://coinranking.com/coin/V8GxkwWow+kaspa-kas'; D='0.116'; E='  -0.51%  ' }
    @{ Row=44; D='28.35'; E='  +8.12%  ' }
    @{ Row=45; D='0.259'; E='  +7.40%  ' }
    @{ Row=47; D='2.10'; E='  +8.87%  ' }
    @{ Row=48; E='  +4.35%  ' }
    @{ Row=49; D='118.71'; E='  +3.34%  ' }
    @{ Row=50; D='0.0₃0516'; E='  +7.09%  ' }
    @{ Row=51; D='2.07'; E='  +2.84%  ' }
)

foreach ($row in $rows) {
    if ($row.ContainsKey("B")) { $ws.Cells.Item($row.Row, 2).Value = $row.B }
    if ($row.ContainsKey("C")) { $ws.Cells.Item($row.Row, 3).Value = $row.C }
    if ($row.ContainsKey("D")) {
        $ws.Cells.Item($row.Row, 4).Value = "'" + $row.D
        $ws.Cells.Item($row.Row, 4).Style = "Normal"
    }
    if ($row.ContainsKey("E")) { $ws.Cells.Item($row.Row, 5).Value = $row.E }
}
